{"js": "// Replace the division-problem text in each table cell with the new\n// problem text, per the commit diff. Each \"from\" string is unique in the\n// document, so an exact (case-sensitive, non-wildcard) search-and-replace\n// on each pair is safe and order-independent.\nconst replacements = [\n  [\"68\u00f76=\", \"28\u00f77=\"],\n  [\"95\u00f74=\", \"89\u00f73=\"],\n  [\"84\u00f77=\", \"54\u00f79=\"],\n  [\"21\u00f73=\", \"95\u00f77=\"],\n  [\"17\u00f75=\", \"86\u00f72=\"],\n  [\"84\u00f76=\", \"16\u00f72=\"],\n  [\"91\u00f72=\", \"49\u00f75=\"],\n  [\"14\u00f73=\", \"31\u00f76=\"],\n  [\"89\u00f72=\", \"69\u00f77=\"],\n  [\"44\u00f79=\", \"41\u00f72=\"],\n  [\"71\u00f74=\", \"74\u00f72=\"],\n  [\"10\u00f73=\", \"36\u00f79=\"],\n  [\"15\u00f76=\", \"14\u00f74=\"],\n  [\"20\u00f78=\", \"97\u00f72=\"],\n  [\"69\u00f72=\", \"42\u00f74=\"],\n  [\"40\u00f79=\", \"61\u00f77=\"],\n  [\"40\u00f78=\", \"78\u00f72=\"],\n  [\"19\u00f72=\", \"68\u00f75=\"],\n  [\"78\u00f79=\", \"49\u00f75=\"],\n  [\"82\u00f78=\", \"76\u00f74=\"],\n  [\"54\u00f79=\", \"77\u00f78=\"],\n  [\"10\u00f76=\", \"93\u00f74=\"],\n  [\"24\u00f75=\", \"49\u00f74=\"],\n  [\"17\u00f78=\", \"73\u00f73=\"],\n  [\"87\u00f79=\", \"14\u00f75=\"],\n];\n\nconst body = context.document.body;\n\n// First, resolve every \"from\" range against the ORIGINAL (unmodified)\n// document. Doing all searches before any edits avoids a later search\n// accidentally matching text that an earlier replacement just inserted\n// (this happens here because one replacement's new text equals another\n// replacement's old text: 84\u00f77= -> 54\u00f79=, and separately 54\u00f79= -> 77\u00f78=).\nconst pending = [];\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  pending.push({ results, to, from });\n}\nawait context.sync();\n\nfor (const { results, to, from } of pending) {\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${from}\", found ${results.items.length}`\n    );\n  }\n  // Each search term is unique in the document; replace the single match.\n  results.items[0].insertText(to, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem text in each table cell, per the commit\n# diff. The data-bearing rows of the (single) table are rows 1, 5, 9, 13,\n# 17 (1-based); the rows in between are spacer/blank rows. Writing\n# directly to each Cell.Range.Text by (row, column) position - rather\n# than a text-based Find/Replace - sidesteps any ambiguity from the fact\n# that some new values coincide with other cells' old values (e.g.\n# 84\u00f77= -> 54\u00f79=, while a different cell holds 54\u00f79= -> 77\u00f78=).\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nif ($table.Rows.Count -lt 17) {\n  throw \"Expected the problem table to have at least 17 rows, found $($table.Rows.Count)\"\n}\nif ($table.Columns.Count -lt 5) {\n  throw \"Expected the problem table to have at least 5 columns, found $($table.Columns.Count)\"\n}\n\n$grid = @{\n  \"1,1\" = \"28\u00f77=\"; \"1,2\" = \"89\u00f73=\"; \"1,3\" = \"54\u00f79=\"; \"1,4\" = \"95\u00f77=\"; \"1,5\" = \"86\u00f72=\";\n  \"5,1\" = \"16\u00f72=\"; \"5,2\" = \"49\u00f75=\"; \"5,3\" = \"31\u00f76=\"; \"5,4\" = \"69\u00f77=\"; \"5,5\" = \"41\u00f72=\";\n  \"9,1\" = \"74\u00f72=\"; \"9,2\" = \"36\u00f79=\"; \"9,3\" = \"14\u00f74=\"; \"9,4\" = \"97\u00f72=\"; \"9,5\" = \"42\u00f74=\";\n  \"13,1\" = \"61\u00f77=\"; \"13,2\" = \"78\u00f72=\"; \"13,3\" = \"68\u00f75=\"; \"13,4\" = \"49\u00f75=\"; \"13,5\" = \"76\u00f74=\";\n  \"17,1\" = \"77\u00f78=\"; \"17,2\" = \"93\u00f74=\"; \"17,3\" = \"49\u00f74=\"; \"17,4\" = \"73\u00f73=\"; \"17,5\" = \"14\u00f75=\";\n}\n\n$rows = @(1, 5, 9, 13, 17)\nforeach ($r in $rows) {\n  for ($c = 1; $c -le 5; $c++) {\n    $key = \"$r,$c\"\n    $cell = $table.Cell($r, $c)\n    $cell.Range.Text = $grid[$key]\n  }\n}\n"}
